# Refactor of the naturalization analysis table:
#  - The analysis for process 790357 done at 11:26:29 (old row 63) is superseded /
#    removed; the newer re-analysis for the same process (old row 65, done at
#    11:37:38) stays.
#  - Rows shift up to fill the gap (old row 64 -> row 63, old row 65 -> row 64).
#  - A new record for process 705567 is appended as the new row 65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the obsolete row 63 (process 790357, first/older analysis).
#    This shifts every row below it up by one, so the former row 64
#    (process 792137) becomes row 63, and the former row 65 (process 790357,
#    newer analysis) becomes row 64 - carrying along all of their cell
#    contents (including the long "T" narrative text) automatically.
$ws.Rows.Item(63).Delete()

# 2) The row that is now row 64 (old row 65) stored its process numbers as
#    text and had "N/A" in column C. Normalize it to match the numeric
#    A/B style used elsewhere in the sheet and clear column C.
$ws.Range("A64").Value = 790357
$ws.Range("B64").Value = 790357
$ws.Range("C64").Value = ""

# 3) Append the new row 65 for process 705567.
#    Force text formatting on the columns that would otherwise be
#    auto-converted by Excel to numbers/percentages (process numbers and
#    the percentage column), matching how they are stored elsewhere in
#    this sheet.
$ws.Range("A65").NumberFormat = "@"
$ws.Range("B65").NumberFormat = "@"
$ws.Range("C65").NumberFormat = "@"
$ws.Range("Q65").NumberFormat = "@"

$ws.Range("A65").Value = "705567"
$ws.Range("B65").Value = "705567"
$ws.Range("C65").Value = "N/A"
$ws.Range("D65").Value = "11 de Mar de 2025"
$ws.Range("E65").Value = "Naturalização Ordinária"
$ws.Range("F65").Value = "Indeferimento"
$ws.Range("G65").Value = "Art. 65, inciso II da Lei nº 13.445/2017; Art. 65, inciso III da Lei nº 13.445/2017; Art. 65, inciso IV da Lei nº 13.445/2017; Não anexou item 8; Não anexou item 4; Não anexou item 3; Não anexou item 2"
$ws.Range("H65").Value = "Indeferimento"
$ws.Range("I65").Value = "Nenhum"
$ws.Range("J65").Value = "Processo indeferido por não atender aos requisitos"
$ws.Range("K65").Value = "✅ ATENDIDO"
$ws.Range("L65").Value = "❌ NÃO ATENDIDO - Prazo de residência não localizado nos campos do sistema"
$ws.Range("M65").Value = "❌ NÃO ATENDIDO - Não anexou item 13 - Comprovante de comunicação em português"
$ws.Range("N65").Value = "❌ NÃO ATENDIDO - Antecedentes criminais inválidos ou não anexados"
$ws.Range("O65").Value = "✅ 0% (0/4)"
$ws.Range("P65").Value = "1/8"
$ws.Range("Q65").Value = "12.5%"
$ws.Range("R65").Value = "17/11/2025"
$ws.Range("S65").Value = "13:35:18"
$ws.Range("T65").Value = "1. Nos termos da legislação, realizadas as diligências necessárias à instrução do presente pedido de Transformação de Naturalização Provisória em Definitiva apresento o presente Relatório Opinativo.`n2. A relação de documentos exigidos pela legislação, foi apresentada integralmente conforme documentos juntados ao processo.`n3. De acordo com a documentação apresentada, o naturalizando possui capacidade civil.`n4. Em relação às condições exigidas para a transformação de naturalização provisória em..."
$ws.Range("U65").Value = "Não atendeu 3 requisito(s)"

# The multi-line text entered above makes Excel apply an automatic custom
# row height; re-run AutoFit so row 65 keeps the sheet's default height,
# consistent with every other row.
$ws.Rows.Item(65).AutoFit()
